$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3 C/D values (session end-time correction)
$ws.Range("C3").Value = 45688.88619232639
$ws.Range("D3").Value = 45688.88624818287

# Append new analytics tracking rows (4-21)
$ws.Range("A4").Value = "MAT111"
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 45688.8971716551
$ws.Range("C4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D4").Value = 45688.89718125
$ws.Range("D4").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A5").Value = "MAT111"
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 45688.89738734953
$ws.Range("C5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D5").Value = 45688.897398125
$ws.Range("D5").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A6").Value = "MAT111"
$ws.Range("B6").Value = 70
$ws.Range("C6").Value = 45688.89738734953
$ws.Range("C6").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D6").Value = 45688.89943899305
$ws.Range("D6").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A7").Value = "MAT111"
$ws.Range("B7").Value = 35
$ws.Range("C7").Value = 45688.89926336805
$ws.Range("C7").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D7").Value = 45688.90001787037
$ws.Range("D7").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A8").Value = "MAT111"
$ws.Range("B8").Value = 180
$ws.Range("C8").Value = 45688.89992260416
$ws.Range("C8").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D8").Value = 45688.93094452546
$ws.Range("D8").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A9").Value = "MAT111"
$ws.Range("B9").Value = 24
$ws.Range("C9").Value = 45688.93066046296
$ws.Range("C9").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D9").Value = 45688.94555541666
$ws.Range("D9").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A10").Value = "MAT111"
$ws.Range("B10").Value = 6
$ws.Range("C10").Value = 45688.94548275463
$ws.Range("C10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D10").Value = 45688.94555541666
$ws.Range("D10").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A11").Value = "MAT141"
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 45688.94568693287
$ws.Range("C11").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D11").Value = 45688.94571783565
$ws.Range("D11").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A12").Value = "MAT141"
$ws.Range("B12").Value = 2
$ws.Range("C12").Value = 45688.94568693287
$ws.Range("C12").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D12").Value = 45688.94579412037
$ws.Range("D12").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A13").Value = "MAT141"
$ws.Range("B13").Value = 2
$ws.Range("C13").Value = 45688.94576165509
$ws.Range("C13").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D13").Value = 45688.94579412037
$ws.Range("D13").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A14").Value = "MAT111"
$ws.Range("B14").Value = -4
$ws.Range("C14").Value = 45688.94782407407
$ws.Range("C14").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D14").Value = 45688.94783728009
$ws.Range("D14").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A15").Value = "MAT111"
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = 45688.94782407407
$ws.Range("C15").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D15").Value = 45688.94946922454
$ws.Range("D15").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A16").Value = "MAT111"
$ws.Range("B16").Value = 4
$ws.Range("C16").Value = 45688.94940972222
$ws.Range("C16").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D16").Value = 45688.94946922454
$ws.Range("D16").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A17").Value = "MAT141"
$ws.Range("B17").Value = 3
$ws.Range("C17").Value = 45688.94949074074
$ws.Range("C17").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D17").Value = 45688.94957386574
$ws.Range("D17").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A18").Value = "MAT141"
$ws.Range("B18").Value = 6
$ws.Range("C18").Value = 45688.94949074074
$ws.Range("C18").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D18").Value = 45688.94986707176
$ws.Range("D18").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A19").Value = "MAT141"
$ws.Range("B19").Value = 19
$ws.Range("C19").Value = 45688.94960648148
$ws.Range("C19").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D19").Value = 45688.94986707176
$ws.Range("D19").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A20").Value = "MAT111"
$ws.Range("B20").Value = 5
$ws.Range("C20").Value = 45688.95012731481
$ws.Range("C20").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D20").Value = 45688.95017834491
$ws.Range("D20").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A21").Value = "MAT111"
$ws.Range("B21").Value = 7
$ws.Range("C21").Value = 45688.95012731481
$ws.Range("C21").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D21").Value = 45688.95017834548
$ws.Range("D21").NumberFormat = "YYYY-MM-DD HH:MM:SS"
